# Apply "small gui tweaks" described in the commit:
#  - Typography sheet: add a new column header "Widget Wildcard Characters" in J3
#  - Translation sheet: add a new text-id row (row 6) "SingleUseId2" / "Disco DCC LCC"

$wb = $excel.ActiveWorkbook

# --- Typography sheet ---------------------------------------------------
$typo = $wb.Worksheets.Item("Typography")
$typo.Range("J3").Value = "Widget Wildcard Characters"
# J4 stays blank, but gets "touched" (present as an empty cell) just like
# its row-4 neighbours (F4/G4/I4) already are in the original sheet.
$typo.Range("J4").NumberFormat = "General"
$typo.Range("J4").Style = "Normal"

# --- Translation sheet ---------------------------------------------------
$trans = $wb.Worksheets.Item("Translation")
$trans.Range("B6").Value = "SingleUseId2"
$trans.Range("C6").Value = "Default"
$trans.Range("D6").Value = "Left"
$trans.Range("E6").Value = "Disco DCC LCC"
$trans.Range("F6").Value = "LTR"
$trans.Range("B6:F6").Style = "Normal"
